$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from
# 45182 (2023-09-13) to 45184 (2023-09-15) for every data row (2..310).
$ws.Range("C2:C310").Value = 45184
